$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update current page (bookmark) for "Researching Information Systems and Computing"
$ws.Range("C11").Value = 202

# Update the active cell selection as recorded in the saved view state
$ws.Range("C19").Select()
